$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3 (rotated from old row 4) ---
$ws.Range("A3").Value = 130822176
$ws.Range("B3").Value = 91804
$ws.Range("E3").Value = 1108
$ws.Range("F3").Value = "Harticka"
$ws.Range("G3").Value = "Pelloporus leporinus"
$ws.Range("H3").Value = "(Fr.) Krieglst."
$ws.Range("Q3").Value = 423906
$ws.Range("R3").Value = 7049117
$ws.Range("AC3").ClearContents()

# --- Row 4 (rotated from old row 3) ---
$ws.Range("A4").Value = 130822165
$ws.Range("AC4").Value = "Ringhack äldre"
$ws.Range("B4").Value = 57884
$ws.Range("E4").Value = 100109
$ws.Range("F4").Value = "Tretåig hackspett"
$ws.Range("G4").Value = "Picoides tridactylus"
$ws.Range("H4").Value = "(Linnaeus, 1758)"
$ws.Range("Q4").Value = 423776
$ws.Range("R4").Value = 7049067

# --- Row 10 (rotated from old row 11) ---
$ws.Range("A10").Value = 130823153
$ws.Range("AC10").Value = "Ringhack, äldre, längs några meter på en granstam. I flerskiktad granskog."
$ws.Range("AI10").Value = "Flerskiktad grandominerad skog med inslag av björk."
$ws.Range("Q10").Value = 423817
$ws.Range("R10").Value = 7049177

# --- Row 11 (rotated from old row 10) ---
$ws.Range("A11").Value = 130823155
$ws.Range("AC11").Value = "Ringhack, äldre, på gran."
$ws.Range("Q11").Value = 423787
$ws.Range("R11").Value = 7049107
$ws.Range("AI11").ClearContents()

# --- Row 19 (rotated from old row 20) ---
$ws.Range("A19").Value = 130822180
$ws.Range("AW19").Value = "Benny Öwre"
$ws.Range("AX19").Value = "Benny Öwre"
$ws.Range("B19").Value = 91804
$ws.Range("E19").Value = 1108
$ws.Range("F19").Value = "Harticka"
$ws.Range("G19").Value = "Pelloporus leporinus"
$ws.Range("H19").Value = "(Fr.) Krieglst."
$ws.Range("P19").Value = "Djupsjö ö, Jmt"
$ws.Range("Q19").Value = 423732
$ws.Range("R19").Value = 7049150
$ws.Range("AC19").ClearContents()
$ws.Range("AF19").ClearContents()
$ws.Range("AH19").ClearContents()
$ws.Range("AJ19").ClearContents()
$ws.Range("AK19").ClearContents()
$ws.Range("AM19").ClearContents()
$ws.Range("AO19").ClearContents()
$ws.Range("J19").ClearContents()
$ws.Range("K19").ClearContents()
$ws.Range("N19").ClearContents()

# --- Row 20 (rotated from old row 21) ---
$ws.Range("A20").Value = 130822182
$ws.Range("Q20").Value = 423748
$ws.Range("R20").Value = 7049159

# --- Row 21 (rotated from old row 19) ---
$ws.Range("A21").Value = 130823158
$ws.Range("AC21").Value = "Små fruktkroppar i en ca 3 meters granhögstubbe."
$ws.Range("AF21").Value = ""
$ws.Range("AH21").Value = "Granskog"
$ws.Range("AJ21").Value = "gran"
$ws.Range("AK21").Value = "Picea abies"
$ws.Range("AM21").Value = "Stående död trädstam/högstubbe"
$ws.Range("AO21").Value = "Standing dead tree/snags # Picea abies"
$ws.Range("AW21").Value = "Kristian Zackrisson"
$ws.Range("AX21").Value = "Kristian Zackrisson"
$ws.Range("B21").Value = 91828
$ws.Range("E21").Value = 5432
$ws.Range("F21").Value = "Granticka"
$ws.Range("G21").Value = "Porodaedalea chrysoloma s.lat."
$ws.Range("H21").Value = ""
$ws.Range("J21").Value = ""
$ws.Range("K21").Value = ""
$ws.Range("N21").Value = ""
$ws.Range("P21").Value = "Svartnäset Djupsjö, Jmt"
$ws.Range("Q21").Value = 423730
$ws.Range("R21").Value = 7049088

# --- Row 32 (rotated from old row 33) ---
$ws.Range("A32").Value = 130822159
$ws.Range("AC32").Value = "Ringhack färska"
$ws.Range("AW32").Value = "Benny Öwre"
$ws.Range("AX32").Value = "Benny Öwre"
$ws.Range("P32").Value = "Djupsjö ö, Jmt"
$ws.Range("Q32").Value = 423805
$ws.Range("R32").Value = 7049145
$ws.Range("AH32").ClearContents()
$ws.Range("AJ32").ClearContents()
$ws.Range("AK32").ClearContents()
$ws.Range("AM32").ClearContents()
$ws.Range("AO32").ClearContents()
$ws.Range("K32").ClearContents()
$ws.Range("L32").ClearContents()
$ws.Range("M32").ClearContents()
$ws.Range("N32").ClearContents()

# --- Row 33 (rotated from old row 34) ---
$ws.Range("A33").Value = 130822158
$ws.Range("AC33").Value = "Ringhack äldre"
$ws.Range("Q33").Value = 423797
$ws.Range("R33").Value = 7049151

# --- Row 34 (rotated from old row 35) ---
$ws.Range("A34").Value = 130822154
$ws.Range("Q34").Value = 423964
$ws.Range("R34").Value = 7049137

# --- Row 35 (rotated from old row 32) ---
$ws.Range("A35").Value = 130823152
$ws.Range("AC35").Value = "Ringhack, äldre, enstaka några meter upp på en gran."
$ws.Range("AH35").Value = "Granskog"
$ws.Range("AJ35").Value = "gran"
$ws.Range("AK35").Value = "Picea abies"
$ws.Range("AM35").Value = "Trädstam på levande träd"
$ws.Range("AO35").Value = "Stem on living tree # Picea abies"
$ws.Range("AW35").Value = "Kristian Zackrisson"
$ws.Range("AX35").Value = "Kristian Zackrisson"
$ws.Range("K35").Value = ""
$ws.Range("L35").Value = ""
$ws.Range("M35").Value = "äldre spår"
$ws.Range("N35").Value = ""
$ws.Range("P35").Value = "Svartnäset Djupsjö, Jmt"
$ws.Range("Q35").Value = 423916
$ws.Range("R35").Value = 7049107
